$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Valor Mora" total
$ws.Range("E11").Value = 156000

# Update "Cant. Periodos" count (7 -> 3)
$ws.Range("F13").Value = 3

# Update the first data row: period 2507 -> 2501, value 48533 -> 52000
$ws.Range("E16").Value = "2501"
$ws.Range("F16").Value = 52000

# Update the second data row: period 2506 -> 2502 (value already 52000)
$ws.Range("E17").Value = "2502"

# Update the last surviving data row (currently row 22): period 2501 -> 2503
$ws.Range("E22").Value = "2503"

# Remove the now-superseded middle rows (old periods 2505, 2504, 2503, 2502)
$ws.Rows("18:21").Delete()
